# Generate Report for Handoff
#
# The file "271fd774-4036-4853-9517-8c7899c337d2" has finished its
# handback cycle ("Handed back: in sync with en-US") and a new handoff
# cycle has begun for it ("Ready for handoff"), while
# "e2b17401-f77b-40b5-a932-0e58c15c719b" now holds the row previously
# occupied by "271fd774..." (status "Handed back: in sync with en-US").
# In effect rows 2 and 3 swap which source file they describe, and the
# "271fd774..." row picks up a new status + refreshed handoff
# timestamp(s).

$wb = $excel.ActiveWorkbook

$OLD = "271fd774-4036-4853-9517-8c7899c337d2"
$NEW = "e2b17401-f77b-40b5-a932-0e58c15c719b"

$OLD_MD_URL = "https://github.com/OpenLocalizationTest/oltest/blob/d0b36118cac31ffca23748382dfbade4852351f3/e2e/$OLD.md"
$NEW_MD_URL = "https://github.com/OpenLocalizationTest/oltest/blob/d0b36118cac31ffca23748382dfbade4852351f3/e2e/$NEW.md"

function Set-RowHyperlink($ws, $cellAddr, $url, $text) {
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $url, [System.Reflection.Missing]::Value, [System.Reflection.Missing]::Value, $text) | Out-Null
}

# ---------------------------------------------------------------
# Overview sheet
# ---------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

$ov.Range("B2").Value2 = "Handed back: in sync with en-US"
$ov.Range("C2").Value2 = "Handed back: in sync with en-US"
$ov.Range("D2").Value2 = "2016-40-12 22:40:21"

$ov.Range("B3").Value2 = "Ready for handoff"
$ov.Range("C3").Value2 = "Ready for handoff"
$ov.Range("D3").Value2 = "2016-42-12 22:42:07"

$ov.Range("A2").Value2 = "$NEW.md"
$ov.Range("A3").Value2 = "$OLD.md"

$ov.Hyperlinks.Delete()
Set-RowHyperlink $ov "A2" $NEW_MD_URL "$NEW.md"
Set-RowHyperlink $ov "A3" $OLD_MD_URL "$OLD.md"

# ---------------------------------------------------------------
# Per-locale sheets (zh-cn, de-de)
# ---------------------------------------------------------------
$locales = @{
    "zh-cn" = @{
        HandoffRow2 = "2016-03-12 22:39:23"
        HandbackRow2 = "2016-03-12 22:41:36"
        HandoffRow3New = "2016-03-12 22:42:03"
        HandbackRow3 = "2016-03-12 22:41:36"
        XlfOldHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e13f1eb4a7b13724543a46158d49bcf93a97c606/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$OLD.99d13eb0991974279a9c4b1b9eb4e9fc5d7b094d.zh-cn.xlf"
        XlfOldHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4a8253f1344d5cc0be424734613a138fd7c97544/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$OLD.99d13eb0991974279a9c4b1b9eb4e9fc5d7b094d.zh-cn.xlf"
        XlfNewHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e13f1eb4a7b13724543a46158d49bcf93a97c606/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$NEW.b0767920fe73b92ca745c2e826ef4d236011b364.zh-cn.xlf"
        XlfNewHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/4a8253f1344d5cc0be424734613a138fd7c97544/ol-handback/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$NEW.b0767920fe73b92ca745c2e826ef4d236011b364.zh-cn.xlf"
        MdOldUrl = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/66064f95239e839506c64a168349c030853ff25d/e2e/$OLD.md"
        MdNewUrl = "https://github.com/OpenLocalizationTestOrg/oltest.zh-cn/blob/66064f95239e839506c64a168349c030853ff25d/e2e/$NEW.md"
        OldXlf = "$OLD.99d13eb0991974279a9c4b1b9eb4e9fc5d7b094d.zh-cn.xlf"
        NewXlf = "$NEW.b0767920fe73b92ca745c2e826ef4d236011b364.zh-cn.xlf"
    }
    "de-de" = @{
        HandoffRow2 = "2016-03-12 22:40:21"
        HandbackRow2 = "2016-03-12 22:41:42"
        HandoffRow3New = "2016-03-12 22:42:07"
        HandbackRow3 = "2016-03-12 22:41:42"
        XlfOldHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9919fb50b2626c72b63ac3706243aaae59f91da1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$OLD.99d13eb0991974279a9c4b1b9eb4e9fc5d7b094d.de-de.xlf"
        XlfOldHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a1b7eccd7cabbe439f14817a122f7c72a6fe11d3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$OLD.99d13eb0991974279a9c4b1b9eb4e9fc5d7b094d.de-de.xlf"
        XlfNewHandoffUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/9919fb50b2626c72b63ac3706243aaae59f91da1/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$NEW.b0767920fe73b92ca745c2e826ef4d236011b364.de-de.xlf"
        XlfNewHandbackUrl = "https://github.com/OpenLocalizationTestOrg/olhandback/blob/a1b7eccd7cabbe439f14817a122f7c72a6fe11d3/ol-handback/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$NEW.b0767920fe73b92ca745c2e826ef4d236011b364.de-de.xlf"
        MdOldUrl = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c487060928825d4976f2d1264898b391c00755ee/e2e/$OLD.md"
        MdNewUrl = "https://github.com/OpenLocalizationTestOrg/oltest.de-de/blob/c487060928825d4976f2d1264898b391c00755ee/e2e/$NEW.md"
        OldXlf = "$OLD.99d13eb0991974279a9c4b1b9eb4e9fc5d7b094d.de-de.xlf"
        NewXlf = "$NEW.b0767920fe73b92ca745c2e826ef4d236011b364.de-de.xlf"
    }
}

foreach ($localeName in @("zh-cn", "de-de")) {
    $cfg = $locales[$localeName]
    $ws = $wb.Worksheets.Item($localeName)

    # Row 2 now describes the $NEW file; status stays "Handed back".
    $ws.Range("A2").Value2 = "$NEW.md"
    $ws.Range("B2").Value2 = ".md"
    $ws.Range("C2").Value2 = "Handed back: in sync with en-US"
    $ws.Range("D2").Value2 = $cfg.NewXlf
    $ws.Range("E2").Value2 = $cfg.HandoffRow2
    $ws.Range("F2").Value2 = "$NEW.md"
    $ws.Range("G2").Value2 = $cfg.NewXlf
    $ws.Range("H2").Value2 = $cfg.HandbackRow2
    $ws.Range("I2").Value2 = "Include"

    # Row 3 now describes the $OLD file; status becomes "Ready for handoff"
    # with a refreshed handoff datetime.
    $ws.Range("A3").Value2 = "$OLD.md"
    $ws.Range("B3").Value2 = ".md"
    $ws.Range("C3").Value2 = "Ready for handoff"
    $ws.Range("D3").Value2 = $cfg.OldXlf
    $ws.Range("E3").Value2 = $cfg.HandoffRow3New
    $ws.Range("F3").Value2 = "$OLD.md"
    $ws.Range("G3").Value2 = $cfg.OldXlf
    $ws.Range("H3").Value2 = $cfg.HandbackRow3
    $ws.Range("I3").Value2 = "Include"

    $ws.Hyperlinks.Delete()
    Set-RowHyperlink $ws "A2" $NEW_MD_URL "$NEW.md"
    Set-RowHyperlink $ws "B2" $NEW_MD_URL ".md"
    Set-RowHyperlink $ws "D2" $cfg.XlfNewHandoffUrl $cfg.NewXlf
    Set-RowHyperlink $ws "F2" $cfg.MdNewUrl "$NEW.md"
    Set-RowHyperlink $ws "G2" $cfg.XlfNewHandbackUrl $cfg.NewXlf

    Set-RowHyperlink $ws "A3" $OLD_MD_URL "$OLD.md"
    Set-RowHyperlink $ws "B3" $OLD_MD_URL ".md"
    Set-RowHyperlink $ws "D3" $cfg.XlfOldHandoffUrl $cfg.OldXlf
    Set-RowHyperlink $ws "F3" $cfg.MdOldUrl "$OLD.md"
    Set-RowHyperlink $ws "G3" $cfg.XlfOldHandbackUrl $cfg.OldXlf
}

Write-Output "done"
